# Load file's cached value into XLCell.CachedValue
# - Add a new cached/formatted string value ("26:31:45") for the
#   elapsed-time ([h]:mm:ss) formatted representation of the TimeSpan in
#   row 7, and point the "GetFormattedString()" column (G) at it while the
#   "GetString()" column (F) keeps pointing at the original "1.02:31:45".
# - Column C was narrowed slightly as a side effect of the refactor.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cell Values")

# GetFormattedString() for the TimeSpan row now reflects the cached,
# number-format-applied value ("26:31:45") rather than the TimeSpan's own
# ToString() ("1.02:31:45") which GetString() still returns in F7.
$ws.Range("G7").Value = "26:31:45"

# Column C width tweak that came along with this change.
$ws.Columns.Item(3).ColumnWidth = 8.99
